# Optuna Attempt (go back with original)
# Updates forecast values on "Forecast Comparison" sheet and summary
# statistics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 56
$ws1.Range("H2").Value = 12.81
$ws1.Range("L2").Value = 0.86

# Row 3
$ws1.Range("D3").Value = 57
$ws1.Range("H3").Value = 11.77
$ws1.Range("L3").Value = 0.9399999999999999

# Row 4
$ws1.Range("H4").Value = 11.66
$ws1.Range("L4").Value = 1.15

# Row 5
$ws1.Range("H5").Value = 10.6
$ws1.Range("L5").Value = 1.15

# Row 6
$ws1.Range("H6").Value = 9.51
$ws1.Range("L6").Value = 0.83

# Row 7
$ws1.Range("H7").Value = 8.59
$ws1.Range("L7").Value = 1.18

# Row 8
$ws1.Range("H8").Value = 7.82
$ws1.Range("L8").Value = 1.1

# Row 9
$ws1.Range("H9").Value = 6.76
$ws1.Range("L9").Value = 1.15

# Row 10
$ws1.Range("H10").Value = 5.7
$ws1.Range("L10").Value = 0.82

# Row 11
$ws1.Range("D11").Value = 51
$ws1.Range("H11").Value = 4.82
$ws1.Range("L11").Value = 0.93

# Row 12
$ws1.Range("D12").Value = 51
$ws1.Range("H12").Value = 3.82
$ws1.Range("L12").Value = 0.85

# Row 13
$ws1.Range("D13").Value = 51
$ws1.Range("H13").Value = 2.79
$ws1.Range("L13").Value = 1.2

# Row 14
$ws1.Range("D14").Value = 51
$ws1.Range("H14").Value = 1.81
$ws1.Range("L14").Value = 1.15

# Row 15
$ws1.Range("D15").Value = 50
$ws1.Range("H15").Value = 0.82
$ws1.Range("L15").Value = 1.06

# Row 16
$ws1.Range("D16").Value = 50
$ws1.Range("L16").Value = 0.8

# Row 17
$ws1.Range("D17").Value = 50
$ws1.Range("L17").Value = 0.93

# --- Sheet: Summary ---
# Values in column B on this sheet are stored as text (numeric-looking
# strings), so prefix with an apostrophe to force text entry and avoid
# Excel auto-converting them to numeric cells.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'839"
$ws2.Range("B10").Value = "'430"
$ws2.Range("B11").Value = "'220"
$ws2.Range("B12").Value = "'57"
$ws2.Range("B14").Value = "'50"
